# Scheduled-runner style market-data refresh: overwrite the computed
# price/profit columns (H, I, J, K, L, M, N) for a handful of leve rows
# across several Sheets, matching the latest Universalis pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 527.1429000000001
$ws.Range("I18").Value = 458
$ws.Range("K18").Value = 458
$ws.Range("M18").Value = -174

$ws.Range("H93").Value = 29800
$ws.Range("J93").Value = 29800
$ws.Range("L93").Value = 29800
$ws.Range("N93").Value = -34792

$ws.Range("H112").Value = 1129.1538
$ws.Range("J112").Value = 1129.1538
$ws.Range("L112").Value = 3387.4614
$ws.Range("N112").Value = -5603.4614

$ws.Range("H129").Value = 162746.53
$ws.Range("J129").Value = 165406.31
$ws.Range("L129").Value = 496218.93
$ws.Range("N129").Value = -506218.93

$ws.Range("H132").Value = 2519.9268
$ws.Range("I132").Value = 2733.7778
$ws.Range("J132").Value = 980.2
$ws.Range("K132").Value = 8201.3334
$ws.Range("L132").Value = 2940.6
$ws.Range("M132").Value = -5671.3334
$ws.Range("N132").Value = -8000.6

$ws.Range("H135").Value = 13518246
$ws.Range("I135").Value = 453.1613
$ws.Range("K135").Value = 4078.4517
$ws.Range("M135").Value = -1543.4517

$ws.Range("H138").Value = 1937.5358
$ws.Range("J138").Value = 2206.5522
$ws.Range("L138").Value = 6619.6566
$ws.Range("N138").Value = -16899.6566

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1426.2128
$ws.Range("I2").Value = 1356.1578
$ws.Range("K2").Value = 1356.1578
$ws.Range("M2").Value = -1243.1578

$ws.Range("H61").Value = 1972.4103
$ws.Range("I61").Value = 1713.9656
$ws.Range("J61").Value = 2721.9
$ws.Range("K61").Value = 1713.9656
$ws.Range("L61").Value = 2721.9
$ws.Range("M61").Value = -1501.9656
$ws.Range("N61").Value = -3145.9

$ws.Range("H97").Value = 861.35895
$ws.Range("I97").Value = 850.5484
$ws.Range("K97").Value = 850.5484
$ws.Range("M97").Value = -354.5484

$ws.Range("H116").Value = 1426.2128
$ws.Range("I116").Value = 1356.1578
$ws.Range("K116").Value = 1356.1578
$ws.Range("M116").Value = 937.8422

$ws.Range("H132").Value = 9217.375
$ws.Range("I132").Value = 1268.537
$ws.Range("K132").Value = 3805.611
$ws.Range("M132").Value = -1275.611

$ws.Range("H136").Value = 1972.4103
$ws.Range("I136").Value = 1713.9656
$ws.Range("J136").Value = 2721.9
$ws.Range("K136").Value = 5141.8968
$ws.Range("L136").Value = 8165.700000000001
$ws.Range("M136").Value = -2591.8968
$ws.Range("N136").Value = -13265.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1426.2128
$ws.Range("I3").Value = 1356.1578
$ws.Range("K3").Value = 1356.1578
$ws.Range("M3").Value = -1242.1578

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 499.81818

$ws.Range("H22").Value = 346.875
$ws.Range("I22").Value = 219.9
$ws.Range("J22").Value = 558.5
$ws.Range("K22").Value = 219.9
$ws.Range("L22").Value = 558.5
$ws.Range("M22").Value = 130.1
$ws.Range("N22").Value = -1258.5

$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H31").Value = 10010.317
$ws.Range("I31").Value = 18987.445
$ws.Range("K31").Value = 18987.445
$ws.Range("M31").Value = -18692.445

$ws.Range("H34").Value = 10010.317
$ws.Range("I34").Value = 18987.445
$ws.Range("K34").Value = 18987.445
$ws.Range("M34").Value = -18785.445

$ws.Range("H86").Value = 13905874
$ws.Range("I86").Value = 9408.333000000001
$ws.Range("K86").Value = 9408.333000000001
$ws.Range("M86").Value = -8285.333000000001

$ws.Range("H89").Value = 13905874
$ws.Range("I89").Value = 9408.333000000001
$ws.Range("K89").Value = 47041.665
$ws.Range("M89").Value = -41425.665

$ws.Range("H95").Value = 13450
$ws.Range("J95").Value = 13450
$ws.Range("L95").Value = 13450
$ws.Range("N95").Value = -18942

$ws.Range("H113").Value = 499.81818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2505.6
$ws.Range("J75").Value = 3500
$ws.Range("L75").Value = 10500
$ws.Range("N75").Value = -12496

$ws.Range("H78").Value = 2505.6
$ws.Range("J78").Value = 3500
$ws.Range("L78").Value = 31500
$ws.Range("N78").Value = -41484

$ws.Range("H131").Value = 796.08
$ws.Range("J131").Value = 812.96906
$ws.Range("L131").Value = 2438.90718
$ws.Range("N131").Value = -12518.90718

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3769.1428
$ws.Range("I126").Value = 3138.4546
$ws.Range("K126").Value = 9415.363799999999
$ws.Range("M126").Value = -6945.363799999999

$ws.Range("H132").Value = 58975.965
$ws.Range("I132").Value = 59575.277
$ws.Range("J132").Value = 57777.332
$ws.Range("K132").Value = 178725.831
$ws.Range("L132").Value = 173331.996
$ws.Range("M132").Value = -176195.831
$ws.Range("N132").Value = -178391.996

$ws.Range("H136").Value = 15584.154
$ws.Range("J136").Value = 15584.154
$ws.Range("L136").Value = 46752.462
$ws.Range("N136").Value = -51852.462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 30000
$ws.Range("J69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31622

$ws.Range("H72").Value = 30000
$ws.Range("J72").Value = 30000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98112

$ws.Range("H94").Value = 6000
$ws.Range("J94").Value = 6000
$ws.Range("L94").Value = 6000
$ws.Range("N94").Value = -7352

$ws.Range("H104").Value = 35000
$ws.Range("J104").Value = 35000
$ws.Range("L104").Value = 35000
$ws.Range("N104").Value = -41988

$ws.Range("H122").Value = 1785649.5
$ws.Range("I122").Value = 3271074
$ws.Range("J122").Value = 3140
$ws.Range("K122").Value = 9813222
$ws.Range("L122").Value = 9420
$ws.Range("M122").Value = -9810772
$ws.Range("N122").Value = -14320

$ws.Range("H132").Value = 2513.1333
$ws.Range("I132").Value = 1822.1111
$ws.Range("J132").Value = 3549.6667
$ws.Range("K132").Value = 5466.3333
$ws.Range("L132").Value = 10649.0001
$ws.Range("M132").Value = -2936.3333
$ws.Range("N132").Value = -15709.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2424.75
$ws.Range("I122").Value = 2100
$ws.Range("J122").Value = 2749.5
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 8248.5
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -13148.5
